# "Updated file with local changes"
#
# Source data: row 1 (A1) held a leftover styled placeholder value (0, bold
# font + thin border + centered/top alignment); row 2 (A2) held the real
# question payload as a single-line Python dict repr string. The local edit
# removes the placeholder row entirely (shifting the payload up into A1,
# which also sheds the one-off bold/border/alignment style that only the
# placeholder used) and reformats the payload text as pretty-printed JSON.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old placeholder row (A1 = 0, specially styled). This shifts the
# question text up from A2 into A1 with no style applied.
$ws.Rows(1).Delete()

# Replace the single-line Python-dict-repr text with the pretty-printed
# JSON-style rendition (indented, double-quoted keys/strings, null literals).
$ws.Range("A1").Value = "questions = [`n    {`n        `"title`": `"Create a function isValidExpression (expression: str) -&gt bool that determines whether the order of the parentheses (), square brackets [], and curly braces {} in a string is mathematically valid.A valid expression has the following characteristics:Each type of opening bracket is closed by the same type of bracket.Brackets are closed in the correct order.All brackets are part of a matching pair.Solve this problem using a list, replicating stack functionality by adding and/or removing elements from the end of the list only.Example 1Input:\`"{[()]}\`"Output:TrueExplanation:Every open bracket has a corresponding closing bracket in the correct order.Example 2Input:\`"{[(])}\`"Output:FalseExplanation:The order of the square bracket and the parenthesis is mismatched.`",`n        `"ques_type`": null,`n        `"options`": [],`n        `"score`": null`n    }`n]"

# Re-fit the row height now that the cell holds embedded line breaks, so no
# stale/explicit row height sticks around.
$ws.Rows(1).AutoFit()
